# (v2.1.1.9270) Add support for Korean, fix ATCs
#
# Fills in previously-blank ATC codes (column E) for many antimicrobials,
# and fixes a handful of oral/IV DDD (defined daily dose) values in
# columns J-M that had been incorrectly populated or left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E (atc) updates: blank inline string -> text value ---
$ws.Range("E3").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("E28").Value = "J01RA07,QJ01RA07"
$ws.Range("E40").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("E80").Value = "NA"
$ws.Range("E81").Value = "NA"
$ws.Range("E82").Value = "NA"
$ws.Range("E94").Value = "NA"
$ws.Range("E96").Value = "NA"
$ws.Range("E101").Value = "NA"
$ws.Range("E104").Value = "NA"
$ws.Range("E106").Value = "NA"
$ws.Range("E110").Value = "NA"
$ws.Range("E116").Value = "NA"
$ws.Range("E117").Value = "NA"
$ws.Range("E118").Value = "NA"
$ws.Range("E123").Value = "NA"
$ws.Range("E128").Value = "NA"
$ws.Range("E129").Value = "NA"
$ws.Range("E135").Value = "NA"
$ws.Range("E137").Value = "NA"
$ws.Range("E138").Value = "NA"
$ws.Range("E139").Value = "NA"
$ws.Range("E142").Value = "NA"
$ws.Range("E146").Value = "J01RA10,QJ01RA10"
$ws.Range("E147").Value = "J01RA12,QJ01RA12"
$ws.Range("E148").Value = "J01RA11,QJ01RA11"
$ws.Range("E150").Value = "NA"
$ws.Range("E151").Value = "NA"
$ws.Range("E153").Value = "NA"
$ws.Range("E161").Value = "NA"
$ws.Range("E177").Value = "NA"
$ws.Range("E181").Value = "NA"
$ws.Range("E183").Value = "NA"
$ws.Range("E188").Value = "J04AM03,QJ04AM03"
$ws.Range("E191").Value = "NA"
$ws.Range("E194").Value = "NA"
$ws.Range("E195").Value = "NA"
$ws.Range("E204").Value = "NA"
$ws.Range("E206").Value = "NA"
$ws.Range("E216").Value = "NA"
$ws.Range("E220").Value = "NA"
$ws.Range("E223").Value = "NA"
$ws.Range("E226").Value = "J01DH51,QJ01DH51"
$ws.Range("E227").Value = "NA"
$ws.Range("E228").Value = "J01DH56,QJ01DH56"
$ws.Range("E233").Value = "NA"
$ws.Range("E237").Value = "NA"
$ws.Range("E238").Value = "NA"
$ws.Range("E245").Value = "NA"
$ws.Range("E247").Value = "J01RA05,QJ01RA05"
$ws.Range("E249").Value = "NA"
$ws.Range("E252").Value = "NA"
$ws.Range("E257").Value = "NA"
$ws.Range("E260").Value = "NA"
$ws.Range("E262").Value = "NA"
$ws.Range("E263").Value = "J01DH52,QJ01DH52"
$ws.Range("E264").Value = "NA"
$ws.Range("E269").Value = "NA"
$ws.Range("E270").Value = "NA"
$ws.Range("E273").Value = "NA"
$ws.Range("E278").Value = "NA"
$ws.Range("E281").Value = "NA"
$ws.Range("E285").Value = "NA"
$ws.Range("E288").Value = "NA"
$ws.Range("E290").Value = "NA"
$ws.Range("E295").Value = "NA"
$ws.Range("E296").Value = "NA"
$ws.Range("E300").Value = "NA"
$ws.Range("E303").Value = "NA"
$ws.Range("E304").Value = "J01RA14,QJ01RA14"
$ws.Range("E305").Value = "J01RA13,QJ01RA13"
$ws.Range("E306").Value = "NA"
$ws.Range("E310").Value = "J01RA09,QJ01RA09"
$ws.Range("E313").Value = "NA"
$ws.Range("E316").Value = "NA"
$ws.Range("E320").Value = "NA"
$ws.Range("E324").Value = "NA"
$ws.Range("E325").Value = "NA"
$ws.Range("E329").Value = "NA"
$ws.Range("E331").Value = "NA"
$ws.Range("E332").Value = "NA"
$ws.Range("E334").Value = "NA"
$ws.Range("E335").Value = "NA"
$ws.Range("E336").Value = "NA"
$ws.Range("E339").Value = "NA"
$ws.Range("E344").Value = "NA"
$ws.Range("E351").Value = "NA"
$ws.Range("E354").Value = "NA"
$ws.Range("E356").Value = "NA"
$ws.Range("E360").Value = "NA"
$ws.Range("E364").Value = "QJ01FG02"
$ws.Range("E365").Value = "NA"
$ws.Range("E366").Value = "NA"
$ws.Range("E367").Value = "NA"
$ws.Range("E369").Value = "NA"
$ws.Range("E372").Value = "NA"
$ws.Range("E375").Value = "J04AM07,QJ04AM07"
$ws.Range("E376").Value = "J04AM02,QJ04AM02"
$ws.Range("E377").Value = "J04AM06,QJ04AM06"
$ws.Range("E378").Value = "J04AM05,QJ04AM05"
$ws.Range("E382").Value = "NA"
$ws.Range("E383").Value = "NA"
$ws.Range("E392").Value = "NA"
$ws.Range("E394").Value = "C10BA04,QC10BA04"
$ws.Range("E402").Value = "J01RA04,QJ01RA04"
$ws.Range("E405").Value = "NA"
$ws.Range("E406").Value = "J04AM01,QJ04AM01"
$ws.Range("E412").Value = "J01EE06"
$ws.Range("E413").Value = "J01EE02,QJ01EW10,QJ51RE01"
$ws.Range("E416").Value = "J01EE05,QJ01EW03"
$ws.Range("E422").Value = "J01EE07,QJ01EW18"
$ws.Range("E423").Value = "NA"
$ws.Range("E429").Value = "J01EE03"
$ws.Range("E431").Value = "J01EE04"
$ws.Range("E436").Value = "NA"
$ws.Range("E439").Value = "NA"
$ws.Range("E440").Value = "NA"
$ws.Range("E441").Value = "NA"
$ws.Range("E443").Value = "NA"
$ws.Range("E445").Value = "NA"
$ws.Range("E447").Value = "NA"
$ws.Range("E450").Value = "NA"
$ws.Range("E459").Value = "NA"
$ws.Range("E460").Value = "J01RA08,QJ01RA08"
$ws.Range("E461").Value = "NA"
$ws.Range("E462").Value = "NA"
$ws.Range("E465").Value = "J04AM04,QJ04AM04"
$ws.Range("E470").Value = "NA"
$ws.Range("E476").Value = "NA"
$ws.Range("E477").Value = "NA"
$ws.Range("E478").Value = "NA"
$ws.Range("E480").Value = "NA"
$ws.Range("E483").Value = "J01EE01"
$ws.Range("E485").Value = "NA"
$ws.Range("E490").Value = "NA"
$ws.Range("E492").Value = "NA"
$ws.Range("E493").Value = "NA"
$ws.Range("E494").Value = "NA"
$ws.Range("E497").Value = "NA"
$ws.Range("E498").Value = "NA"

# --- Columns J/K (oral_ddd/oral_units) and L/M (iv_ddd/iv_units) updates ---
# Row 11: Amoxicillin/sulbactam - add oral and IV DDD
$ws.Range("J11").Value = 1.5
$ws.Range("K11").Value = "g"
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = "g"

# Row 111: Cefpodoxime/clavulanic acid - add oral DDD
$ws.Range("J111").Value = 0.4
$ws.Range("K111").Value = "g"

# Row 117: Ceftaroline - remove IV DDD
$ws.Range("L117").Value = "#N/A"
$ws.Range("M117").Value = "#N/A"

# Row 120: Ceftazidime/avibactam - add IV DDD
$ws.Range("L120").Value = 6
$ws.Range("M120").Value = "g"

# Row 136: Cefuroxime/metronidazole - add oral DDD
$ws.Range("J136").Value = 0.5
$ws.Range("K136").Value = "g"

# Row 138: Cephradine - remove oral and IV DDD
$ws.Range("J138").Value = "#N/A"
$ws.Range("K138").Value = "#N/A"
$ws.Range("L138").Value = "#N/A"
$ws.Range("M138").Value = "#N/A"

# Row 325: Panipenem - remove IV DDD
$ws.Range("L325").Value = "#N/A"
$ws.Range("M325").Value = "#N/A"

# Row 342: Piperacillin/sulbactam - add IV DDD
$ws.Range("L342").Value = 14
$ws.Range("M342").Value = "g"

# Row 447: Tebipenem - remove oral DDD
$ws.Range("J447").Value = "#N/A"
$ws.Range("K447").Value = "#N/A"
